$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 updates (average_county_temperature, worst_ashp_cop, best_ashp_cop)
$ws.Range("I8").Value = 19.79629629629628
$ws.Range("N8").Value = 1.98600466835246
$ws.Range("O8").Value = 2.18975222777657

# Row 10 updates (average_county_temperature, worst_ashp_cop, best_ashp_cop)
$ws.Range("I10").Value = 15.74228395061728
$ws.Range("N10").Value = 1.904889690449167
$ws.Range("O10").Value = 2.090295475371289
